$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-05 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-06 Thursday", 2) | Out-Null
$d.Content.Find.Execute("125×5=625", $true, $false, $false, $false, $false, $true, 1, $false, "960×7=6720", 2) | Out-Null
$d.Content.Find.Execute("760×8=6080", $true, $false, $false, $false, $false, $true, 1, $false, "765×5=3825", 2) | Out-Null
$d.Content.Find.Execute("788×8=6304", $true, $false, $false, $false, $false, $true, 1, $false, "409×4=1636", 2) | Out-Null
$d.Content.Find.Execute("888×7=6216", $true, $false, $false, $false, $false, $true, 1, $false, "364×2=728", 2) | Out-Null
$d.Content.Find.Execute("102×7=714", $true, $false, $false, $false, $false, $true, 1, $false, "817×2=1634", 2) | Out-Null
$d.Content.Find.Execute("874×6=5244", $true, $false, $false, $false, $false, $true, 1, $false, "574×6=3444", 2) | Out-Null
$d.Content.Find.Execute("887×3=2661", $true, $false, $false, $false, $false, $true, 1, $false, "958×9=8622", 2) | Out-Null
$d.Content.Find.Execute("770×2=1540", $true, $false, $false, $false, $false, $true, 1, $false, "698×2=1396", 2) | Out-Null
$d.Content.Find.Execute("925×6=5550", $true, $false, $false, $false, $false, $true, 1, $false, "579×2=1158", 2) | Out-Null
$d.Content.Find.Execute("300×7=2100", $true, $false, $false, $false, $false, $true, 1, $false, "985×8=7880", 2) | Out-Null
$d.Content.Find.Execute("759×4=3036", $true, $false, $false, $false, $false, $true, 1, $false, "936×5=4680", 2) | Out-Null
$d.Content.Find.Execute("718×5=3590", $true, $false, $false, $false, $false, $true, 1, $false, "445×4=1780", 2) | Out-Null
$d.Content.Find.Execute("467×2=934", $true, $false, $false, $false, $false, $true, 1, $false, "995×5=4975", 2) | Out-Null
$d.Content.Find.Execute("711×2=1422", $true, $false, $false, $false, $false, $true, 1, $false, "491×3=1473", 2) | Out-Null
$d.Content.Find.Execute("738×9=6642", $true, $false, $false, $false, $false, $true, 1, $false, "176×9=1584", 2) | Out-Null
$d.Content.Find.Execute("127×9=1143", $true, $false, $false, $false, $false, $true, 1, $false, "440×2=880", 2) | Out-Null
$d.Content.Find.Execute("317×7=2219", $true, $false, $false, $false, $false, $true, 1, $false, "807×3=2421", 2) | Out-Null
$d.Content.Find.Execute("779×3=2337", $true, $false, $false, $false, $false, $true, 1, $false, "304×6=1824", 2) | Out-Null
$d.Content.Find.Execute("265×9=2385", $true, $false, $false, $false, $false, $true, 1, $false, "111×7=777", 2) | Out-Null
$d.Content.Find.Execute("333×5=1665", $true, $false, $false, $false, $false, $true, 1, $false, "670×8=5360", 2) | Out-Null
$d.Content.Find.Execute("985×7=6895", $true, $false, $false, $false, $false, $true, 1, $false, "473×9=4257", 2) | Out-Null
$d.Content.Find.Execute("359×5=1795", $true, $false, $false, $false, $false, $true, 1, $false, "589×7=4123", 2) | Out-Null
$d.Content.Find.Execute("714×6=4284", $true, $false, $false, $false, $false, $true, 1, $false, "589×6=3534", 2) | Out-Null
$d.Content.Find.Execute("318×8=2544", $true, $false, $false, $false, $false, $true, 1, $false, "494×3=1482", 2) | Out-Null
$d.Content.Find.Execute("404×7=2828", $true, $false, $false, $false, $false, $true, 1, $false, "711×3=2133", 2) | Out-Null
